# Weekly refresh of the fruit/vegetable price logs: the Fecha (date) and
# price-related columns are reshuffled across the existing data rows
# (rows 2-39) of the sheet, as if each row's record was replaced by
# another day's record while keeping the row's fixed descriptive
# columns (Mercado, Region, Categoria, etc.) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: targetRow -> sourceRow. The target row receives the
# "payload" (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Origen, Precio $/Kg) that used to belong to the source row.
$rowMap = @{
    2=11; 3=6; 4=35; 5=3; 6=13; 7=24; 8=7; 9=26; 10=8;
    11=4; 12=29; 13=17; 14=16; 15=27; 16=23; 17=15; 18=14; 19=33; 20=28;
    21=36; 22=2; 23=39; 24=22; 25=30; 26=18; 27=19; 28=5; 29=34; 30=32;
    31=12; 32=10; 33=38; 34=9; 35=37; 36=20; 37=31; 38=25; 39=21
}

# Columns that move together as the payload for each row.
# D=4 (Fecha), I=9 (Calidad), J=10 (Volumen), K=11 (Precio minimo),
# L=12 (Precio maximo), M=13 (Precio promedio ponderado), O=15 (Origen),
# P=16 (Precio $/Kg)
$payloadCols = @(4, 9, 10, 11, 12, 13, 15, 16)

# First snapshot every row's current payload values, since the
# permutation is a single large cycle and later writes must not clobber
# values that are still needed as a source for another row.
$snapshot = @{}
for ($r = 2; $r -le 39; $r++) {
    $rowVals = @{}
    foreach ($c in $payloadCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now apply the new payload to each target row using the snapshotted
# source-row values.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $payloadCols) {
        $ws.Cells.Item($targetRow, $c).Value = $srcVals[$c]
    }
}
